# Fill in Jesse's Week 9 local copies of the Task Summary Sheet and
# Activity Log Summary Sheet with his task data.

$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET -------------------------------------------------
$task = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Name + Week number
$task.Range("C1").Value = "Jesse Hare"
$task.Range("E1").Value = 9

# Task rows
$task.Range("A3").Value = "Project Build"
$task.Range("A4").Value = "Project Build"

$task.Range("B4").Value = "Implement ideas/suggestions from Team meeting"
$task.Range("B3").Value = "Code improvement, optimisation, ironing out bugs"

$task.Range("C3").Value = 5
$task.Range("D3").Value = 10
$task.Range("E3").Value = 0

$task.Range("C4").Value = 10
$task.Range("D4").Value = 10
$task.Range("E4").Value = 0

# --- ACTIVITY LOG SUMMARY SHEET -----------------------------------------
$summary = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$summary.Range("D1").Value = "Jesse Hare"

$summary.Range("A4").Value = "Project Build"
$summary.Range("B4").Value = 16
$summary.Range("C4").Value = 4
